# Edit: add "Sprint Backlog - Meeting 3" and "Sprint BackLog 3" sheets to the
# workbook (before "Team Members"), populate the new Sprint Backlog - Meeting 3
# sheet with the Product Backlog content plus eight new backlog items, and
# update a couple of selections / the active tab to match.

$wb = $excel.ActiveWorkbook

$sprintBacklog2 = $wb.Worksheets.Item("Sprint Backlog 2")

# Create the two new sheets in between "Sprint Backlog 2" and "Team Members".
# Creating "Sprint Backlog - Meeting 3" first gives it the lower sheetId, then
# "Sprint BackLog 3" is added directly after it (both still land before
# "Team Members").
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sprintBacklog2)
$ws.Name = "Sprint Backlog - Meeting 3"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws)
$ws2.Name = "Sprint BackLog 3"

# Populate "Sprint Backlog - Meeting 3" - header row plus the product backlog
# items (same as "Product Backlog - Meeting 2") plus eight new user stories.
$ws.Range("A1").Value = 'Product Backlog'
$ws.Range("B1").Value = 'Bus Value'
$ws.Range("C1").Value = 'Effort'
$ws.Range("D1").Value = 'ROI'
$ws.Range("A2").Value = 'As a payment provider I want an accept Payment API available so that I can provide my services to online stores'
$ws.Range("A3").Value = 'As a payment provider  I want an Payment API available so that I can troubleshoot payment issues'
$ws.Range("A4").Value = 'As a payment provider I want an Payment API available so that I can see all payments on a date'
$ws.Range("A5").Value = 'As a payment provider I want Payment info in a database so that I can see all my data for reports '
$ws.Range("A6").Value = 'As a payment provider I want customer specific portal so that my customers can look at their transaction data directly '
$ws.Range("A7").Value = 'As a hotel customer I want a webpage details page so that I can view my payment status'
$ws.Range("A8").Value = 'As a hotel owner I want a payment details portal so that I can assist my customers'
$ws.Range("A9").Value = 'As a hotel customer I want a user portal so that I can cancel my booking '
$ws.Range("A10").Value = 'As a hotel customer I want a user portal so that I can edit my booking '
$ws.Range("A11").Value = 'As a customer I want to have member booking so that I do not have to remember my contact info'
$ws.Range("A12").Value = 'As a customer I want to have member booking so that I can get special discounts off my booking'
$ws.Range("A13").Value = 'As a owner I want to support loyality points to that customers can redeem for a deal'
$ws.Range("A14").Value = 'As an owner I want to see what my customers are looking at so that I can provide customer specific offers'
$ws.Range("A15").Value = 'As a payment provideer I want a payment API that can hamdle different tax rates'
$ws.Range("A16").Value = 'As an owner I want to review/feedback page visible by owner so that I can improve customer satisfaction '
$ws.Range("A17").Value = 'As a customer I want to get alerts on bookings and new deals so that I can know about deals'
$ws.Range("A18").Value = 'As an owner I want to support A/B testing so that I can optimize ad placement for revenue'
$ws.Range("A19").Value = 'As a customer I want to see my pervious bookings so that I can see what I paid ( compare prices)'
$ws.Range("A20").Value = 'As a customer I want to get an email notification of payment so that I can have proof of payment'
$ws.Range("A21").Value = 'As a hotel owner I want to see the resurant menu so that customers will know what food options there are'
$ws.Range("A22").Value = 'As a customer I want to have an option to book a meeting hall so that I can online book an event'
$ws.Range("A23").Value = 'As a hotel owner I want a special role so that I can have different access for different users ( admin/users)'
$ws.Range("A24").Value = 'As a hotel owner I want to implement a live chat for employees so that I can assist my customers in real time'
$ws.Range("A25").Value = 'As a customer I want to schedue an airport pickup so that I be picked up by the hotel airport shuttle'
$ws.Range("A26").Value = 'As a customer I want to re-schedue an airport pickup so that I be picked up by the hotel airport shuttle'

# Formatting to match the other Product/Sprint Backlog sheets.
$ws.Columns.Item(1).ColumnWidth = 96.16666666666667

# Leave the selection on the empty row just below the data, then zoom + make
# this the active sheet/tab.
$ws.Range("A27").Select() | Out-Null
$excel.ActiveWindow.Zoom = 83

# "Sprint BackLog 3" stays empty; just move the cursor to P31.
$ws2.Range("P31").Select() | Out-Null

# Update the selection on "Product Backlog - Meeting 2".
$productBacklog2 = $wb.Worksheets.Item("Product Backlog - Meeting 2")
$productBacklog2.Range("B1:D1").Select() | Out-Null

# Finally, activate "Sprint Backlog - Meeting 3" so it becomes the selected
# tab (matches the workbook's activeTab moving from "Sprint Backlog 2" to
# this new sheet).
$ws.Activate() | Out-Null
